$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Update the test-data row (row 2) with the new AUTOTEST naming convention
# and the new GL account strings / natural account number.
$ws.Range("I2").Value = "AUTOTEST"
$ws.Range("J2").Value = "AUTOTEST"
$ws.Range("M2").Value = "220.31863.204021.10001.001.00000.000"
$ws.Range("N2").Value = "220.31863.204021.10001.000.00000.000"
$ws.Range("Q2").Value = "AUTOTEST"
$ws.Range("S2").Value = "AUTOTEST"
$ws.Range("T2").Value = "AUTOTEST"
$ws.Range("W2").Value = 204021

# Move the active cell selection on the Input_Value sheet
[void]$ws.Range("G9").Select()
